# Update crypto price/volume figures to the latest scrape (GitHub Actions run).
# Rows 28/29 (Toncoin/Cosmos) also swapped rank position in the source feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.797.74"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "2.543.22"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'303.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").Value = "'97.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.72%  "

$ws.Range("D7").Value = "'0.577"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D9").Value = "'0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("D10").Value = "'36.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("D11").Value = "'0.0822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.27%  "

$ws.Range("E12").Value = "  +0.73%  "

$ws.Range("D13").Value = "'7.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.23%  "

$ws.Range("D14").Value = "2.936.59"
$ws.Range("E14").Value = "  -0.09%  "

$ws.Range("D15").Value = "2.614.81"
$ws.Range("E15").Value = "  +4.03%  "

$ws.Range("D16").Value = "'15.06"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.49%  "

$ws.Range("D17").Value = "'0.865"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").Value = "42.850.22"
$ws.Range("E18").Value = "  +0.35%  "

$ws.Range("E19").Value = "  +3.07%  "

$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").Value = "  +0.06%  "

$ws.Range("D21").Value = "'6.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("D22").Value = "'71.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "

$ws.Range("D23").Value = "'253.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").Value = "'2.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("E25").Value = "  -2.22%  "

$ws.Range("D26").Value = "'27.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.47%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'10.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.76%  "

$ws.Range("D30").Value = "'37.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").Value = "'6.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.01%  "

$ws.Range("D32").Value = "'156.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.80%  "

$ws.Range("D33").Value = "'19.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +13.82%  "

$ws.Range("D34").Value = "'2.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("D35").Value = "'3.31"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.26%  "

$ws.Range("D36").Value = "'0.0795"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.15%  "

$ws.Range("D37").Value = "'2.62"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.74%  "

$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("D39").Value = "'25.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.46%  "

$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").Value = "'2.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +28.58%  "

$ws.Range("D42").Value = "'3.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("D44").Value = "2.100.09"
$ws.Range("E44").Value = "  +0.85%  "

$ws.Range("D45").Value = "'0.0304"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.34%  "

$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").Value = "'86.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.13%  "

$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("D49").Value = "2.794.78"
$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").Value = "'73.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.47%  "

$ws.Range("D51").Value = "'0.191"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.67%  "
